$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "[-, 'MEC-1NA-Elet. Dig. Bas.', -, -]"

$ws.Range("C18").Value = $newValue
$ws.Range("C19").Value = $newValue
$ws.Range("C20").Value = $newValue
$ws.Range("C21").Value = $newValue
